# Adds the "28 febbraio - 5 marzo 2022" weekly monitoring row to each of
# the four sheets (Classi, Alunni in presenza, Alunni, Personale scolastico),
# matching the upstream "Add files via upload" commit.

$wb = $excel.ActiveWorkbook
$weekLabel = "28 febbraio - 5 marzo 2022"

# -----------------------------------------------------------------
# Sheet "Classi" (sheet1): add row 9
# -----------------------------------------------------------------
$wsClassi = $wb.Worksheets.Item("Classi")

$wsClassi.Range("A9").Value = $weekLabel

$wsClassi.Range("B9").Value = 5926
$wsClassi.Range("B9").NumberFormat = "#,##0"
$wsClassi.Range("C9").Value = 8157
$wsClassi.Range("C9").NumberFormat = "#,##0"
$wsClassi.Range("D9").Value = 0.72599999999999998
$wsClassi.Range("D9").NumberFormat = "0.0%"
$wsClassi.Range("E9").Value = 376464
$wsClassi.Range("E9").NumberFormat = "#,##0"
$wsClassi.Range("F9").Value = 273255
$wsClassi.Range("F9").NumberFormat = "#,##0"
$wsClassi.Range("G9").Value = 0.72599999999999998
$wsClassi.Range("G9").NumberFormat = "0.0%"
$wsClassi.Range("H9").Value = 272998
$wsClassi.Range("H9").NumberFormat = "#,##0"
$wsClassi.Range("I9").Value = 18307
$wsClassi.Range("I9").NumberFormat = "#,##0"
$wsClassi.Range("J9").Value = 0.99900000000000011
$wsClassi.Range("J9").NumberFormat = "0.0%"
$wsClassi.Range("K9").Value = 0.067000000000000004
$wsClassi.Range("K9").NumberFormat = "0.0%"
$wsClassi.Range("N9").Value = 257
$wsClassi.Range("O9").Value = 0.001
$wsClassi.Range("O9").NumberFormat = "0.0%"
$wsClassi.Range("O9").Font.Color = 0

# -----------------------------------------------------------------
# Sheet "Alunni in presenza" (sheet2): add row 9
# -----------------------------------------------------------------
$wsAlunniPres = $wb.Worksheets.Item("Alunni in presenza")

$wsAlunniPres.Range("A9").Value = $weekLabel
$wsAlunniPres.Range("B9").Value = 7393168
$wsAlunniPres.Range("B9").NumberFormat = "#,##0"
$wsAlunniPres.Range("C9").Value = 5342389
$wsAlunniPres.Range("C9").NumberFormat = "#,##0"
$wsAlunniPres.Range("D9").Value = 0.72299999999999998
$wsAlunniPres.Range("D9").NumberFormat = "0.0%"
$wsAlunniPres.Range("E9").Value = 5208193
$wsAlunniPres.Range("E9").NumberFormat = "#,##0"
$wsAlunniPres.Range("F9").Value = 0.97499999999999998
$wsAlunniPres.Range("F9").NumberFormat = "0.0%"

# -----------------------------------------------------------------
# Sheet "Alunni" (sheet3): add rows 30, 31, 32 (row 29 stays a gap,
# matching the existing pattern of blank separator rows that are
# entirely omitted from the sheet)
# -----------------------------------------------------------------
$wsAlunni = $wb.Worksheets.Item("Alunni")

$wsAlunni.Range("A30").Value = $weekLabel
$wsAlunni.Range("B30").Value = "Infanzia"
$wsAlunni.Range("C30").Value = 608188
$wsAlunni.Range("C30").NumberFormat = "#,##0"
$wsAlunni.Range("D30").Value = 594246
$wsAlunni.Range("D30").NumberFormat = "#,##0"
$wsAlunni.Range("E30").Value = 13942
$wsAlunni.Range("E30").NumberFormat = "#,##0"
$wsAlunni.Range("F30").Value = 0.023
$wsAlunni.Range("F30").NumberFormat = "0.0%"

$wsAlunni.Range("A31").Value = $weekLabel
$wsAlunni.Range("B31").Value = "Primaria"
$wsAlunni.Range("C31").Value = 1687275
$wsAlunni.Range("C31").NumberFormat = "#,##0"
$wsAlunni.Range("D31").Value = 1643817
$wsAlunni.Range("D31").NumberFormat = "#,##0"
$wsAlunni.Range("E31").Value = 43458
$wsAlunni.Range("E31").NumberFormat = "#,##0"
$wsAlunni.Range("F31").Value = 0.026000000000000002
$wsAlunni.Range("F31").NumberFormat = "0.0%"

$wsAlunni.Range("A32").Value = $weekLabel
$wsAlunni.Range("B32").Value = "Sec. 1° e 2° Grado"
$wsAlunni.Range("C32").Value = 3046926
$wsAlunni.Range("C32").NumberFormat = "#,##0"
$wsAlunni.Range("D32").Value = 2970130
$wsAlunni.Range("D32").NumberFormat = "#,##0"
$wsAlunni.Range("E32").Value = 76796
$wsAlunni.Range("E32").NumberFormat = "#,##0"
$wsAlunni.Range("F32").Value = 0.025
$wsAlunni.Range("F32").NumberFormat = "0.0%"

# -----------------------------------------------------------------
# Sheet "Personale scolastico" (sheet4): add row 9
# -----------------------------------------------------------------
$wsPersonale = $wb.Worksheets.Item("Personale scolastico")

$wsPersonale.Range("A9").Value = $weekLabel
$wsPersonale.Range("B9").Value = 775867
$wsPersonale.Range("B9").NumberFormat = "#,##0"
$wsPersonale.Range("C9").Value = 558524
$wsPersonale.Range("C9").NumberFormat = "#,##0"
$wsPersonale.Range("D9").Value = 0.72
$wsPersonale.Range("D9").NumberFormat = "0.0%"
$wsPersonale.Range("E9").Value = 543784
$wsPersonale.Range("E9").NumberFormat = "#,##0"
$wsPersonale.Range("F9").Value = 0.97400000000000009
$wsPersonale.Range("F9").NumberFormat = "0.0%"
$wsPersonale.Range("G9").Value = 204526
$wsPersonale.Range("G9").NumberFormat = "#,##0"
$wsPersonale.Range("H9").Value = 148016
$wsPersonale.Range("H9").NumberFormat = "#,##0"
$wsPersonale.Range("I9").Value = 0.72400000000000009
$wsPersonale.Range("I9").NumberFormat = "0.0%"
$wsPersonale.Range("J9").Value = 145068
$wsPersonale.Range("J9").NumberFormat = "#,##0"
$wsPersonale.Range("K9").Value = 0.98
$wsPersonale.Range("K9").NumberFormat = "0.0%"

# -----------------------------------------------------------------
# Restore the on-screen selection / active-sheet state to match the
# author's saved view (Personale scolastico ends up the active tab,
# matching the unchanged activeTab="3" in workbook.xml).
# -----------------------------------------------------------------
$wsClassi.Range("O9").Select()
$wsAlunniPres.Range("G9").Select()
$wsAlunni.Range("G32").Select()
$wsPersonale.Range("B10").Select()
